# Switch focus from PM10 to PM2.5
#
# The source paragraph contains a single run:
#   "...individual health. PM10 which measures air particulates will be considered."
# The target OOXML splits this into three runs:
#   "...individual health. PM" | "2.5" | " which measures air particulates will be considered."
#
# We locate "PM10", replace the "10" portion with "2.5", then force a run
# boundary around the newly inserted "2.5" text (via a transient bookmark,
# which Word uses internally to split runs at its start/end, and which
# leaves no visible trace once removed) so the three segments are written
# out as three separate <w:r> elements with no residual formatting.

$d = $word.ActiveDocument

$find = $d.Content
$find.Find.Execute("PM10")
$matchStart = $find.Start
$matchEnd = $find.End

# "PM" occupies the first two characters of the match; "10" is the rest.
$pmEnd = $matchStart + 2

# Replace "10" with "2.5".
$digits = $d.Range($pmEnd, $matchEnd)
$digits.Text = "2.5"

# The newly-inserted "2.5" now spans pmEnd .. pmEnd+3.
$newValue = $d.Range($pmEnd, $pmEnd + 3)

# Force "2.5" into its own run (distinct from the surrounding text) without
# leaving any formatting residue, by briefly bookmarking it.
$d.Bookmarks.Add("__tmp_split__", $newValue)
$d.Bookmarks("__tmp_split__").Delete()
